$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple is (row, new A value (date serial), new B value)
$data = @(
    @(2, 46074.01041666666, 2302.824),
    @(3, 46074.02083333334, 2290.375),
    @(4, 46074.03125, 2283.566),
    @(5, 46074.04166666666, 2277.197),
    @(6, 46074.05208333334, 2273.998),
    @(7, 46074.0625, 2258.526),
    @(8, 46074.07291666666, 2256.114),
    @(9, 46074.08333333334, 2240.242),
    @(10, 46074.09375, 2232.002),
    @(11, 46074.10416666666, 2228.131),
    @(12, 46074.11458333334, 2224.797),
    @(13, 46074.125, 2220.095),
    @(14, 46074.13541666666, 2222.537),
    @(15, 46074.14583333334, 2214.861),
    @(16, 46074.15625, 2207.065),
    @(17, 46074.16666666666, 2206.985),
    @(18, 46074.17708333334, 2187.612),
    @(19, 46074.1875, 2180.397),
    @(20, 46074.19791666666, 2158.825),
    @(21, 46074.20833333334, 2152.32),
    @(22, 46074.21875, 2122.388),
    @(23, 46074.22916666666, 2111.598),
    @(24, 46074.23958333334, 2108.99),
    @(25, 46074.25, 2098.661),
    @(26, 46074.26041666666, 2060.225),
    @(27, 46074.27083333334, 2060.313),
    @(28, 46074.28125, 2051.123),
    @(29, 46074.29166666666, 2042.162),
    @(30, 46074.30208333334, 2015.771),
    @(31, 46074.3125, 2005.438),
    @(32, 46074.32291666666, 1995.037),
    @(33, 46074.33333333334, 1984.012),
    @(34, 46074.34375, 1953.954),
    @(35, 46074.35416666666, 1938.69),
    @(36, 46074.36458333334, 1929.549),
    @(37, 46074.375, 1916.242),
    @(38, 46074.38541666666, 1898.887),
    @(39, 46074.39583333334, 1888.365),
    @(40, 46074.40625, 1878.002),
    @(41, 46074.41666666666, 1867.168),
    @(42, 46074.42708333334, 1853.442),
    @(43, 46074.4375, 1842.567),
    @(44, 46074.44791666666, 1831.311),
    @(45, 46074.45833333334, 1810.691),
    @(46, 46074.46875, 1795.835),
    @(47, 46074.47916666666, 1778.907),
    @(48, 46074.48958333334, 1760.429),
    @(49, 46074.5, 1733.07),
    @(50, 46074.51041666666, 1699.638),
    @(51, 46074.52083333334, 1670.644),
    @(52, 46074.53125, 1650.522),
    @(53, 46074.54166666666, 1630.916),
    @(54, 46074.55208333334, 1599.797),
    @(55, 46074.5625, 1583.454),
    @(56, 46074.57291666666, 1559.865),
    @(57, 46074.58333333334, 1545.736),
    @(58, 46074.59375, 1508.666),
    @(59, 46074.60416666666, 1489.409),
    @(60, 46074.61458333334, 1478.327),
    @(61, 46074.625, 1458),
    @(62, 46074.63541666666, 1427.209),
    @(63, 46074.64583333334, 1417.418),
    @(64, 46074.65625, 1399.506),
    @(65, 46074.66666666666, 1381.071),
    @(66, 46074.67708333334, 1381.851),
    @(67, 46074.6875, 1363.691),
    @(68, 46074.69791666666, 1383.396),
    @(69, 46074.70833333334, 1364.168),
    @(70, 46074.71875, 1329.668),
    @(71, 46074.72916666666, 1309.73),
    @(72, 46074.73958333334, 1290.437),
    @(73, 46074.75, 1271.469),
    @(74, 46074.76041666666, 1237.685),
    @(75, 46074.77083333334, 1212.919),
    @(76, 46074.78125, 1186.775),
    @(77, 46074.79166666666, 1162.391),
    @(78, 46074.80208333334, 1121.232),
    @(79, 46074.8125, 1096.588),
    @(80, 46074.82291666666, 1072.018),
    @(81, 46074.83333333334, 1047.866),
    @(82, 46074.84375, 1015.289),
    @(83, 46074.85416666666, 994.963),
    @(84, 46074.86458333334, 973.862),
    @(85, 46074.875, 933.394),
    @(86, 46074.88541666666, 903.005),
    @(87, 46074.89583333334, 884.022),
    @(88, 46074.90625, 866.254),
    @(89, 46074.91666666666, 848.0700000000001),
    @(90, 46074.92708333334, 831.002),
    @(91, 46074.9375, 825.068),
    @(92, 46074.94791666666, 817.423),
    @(93, 46074.95833333334, 810.174),
    @(94, 46074.96875, 0),
    @(95, 46074.97916666666, 0),
    @(96, 46074.98958333334, 0),
    @(97, 46075, 0)
)

foreach ($item in $data) {
    $row = $item[0]
    $aVal = $item[1]
    $bVal = $item[2]
    $ws.Cells.Item($row, 1).Value = $aVal
    $ws.Cells.Item($row, 2).Value = $bVal
}
